# Update the "Metadata" worksheet of the CodeSystem workbook:
#  - Version 5.0.0 -> 6.0.0
#  - Date updated to the new publication timestamp
#  - Publisher gets a value ("Alvearie Team")
#  - The (duplicated) "Contact" / "No display for ContactDetail" row is
#    dropped and replaced by a single "Jurisdiction" / "United States of
#    America" row
#  - "Case Sensitive" gets a value of "true"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the first of the two duplicate "Contact" rows (row 10); this
# shifts every row below it up by one, so the second "Contact" row
# (formerly row 11) becomes row 10 and gets turned into "Jurisdiction".
$ws.Rows.Item(10).Delete() | Out-Null

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher
$ws.Range("B9").Value = "Alvearie Team"

# Jurisdiction (previously the second "Contact" row)
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive -> "true"
# NOTE: assigning the literal text "true"/"false" directly to a Range's
# Value/Value2/Formula auto-coerces it to an Excel boolean (TRUE/FALSE),
# same as typing it into a cell in real Excel. To store the literal text
# "true" we compute it as a formula result (never autocoerced) in a scratch
# cell, then copy/paste-special as values into the target cell, which
# writes a plain text value without disturbing the cell's existing style.
$scratch = $ws.Range("Z1")
$scratch.Formula = '=CONCATENATE("tr","ue")'
$scratch.Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4163) | Out-Null
$scratch.EntireColumn.Delete() | Out-Null

Write-Output "done"
